$d = $word.ActiveDocument

# The document currently ends with a single empty paragraph (the last
# paragraph in the story). We turn that paragraph into the new
# "Variables de respuesta" heading, and insert four more paragraphs
# after it (while everything is still in its plain/default state, so
# the new paragraphs don't inherit the heading's direct formatting),
# finishing with a paragraph that holds the closing sentence.

$heading = $d.Paragraphs.Last

# Create the five new (still-blank, still-default-formatted) paragraphs
# first, by repeatedly splitting the last paragraph in the story.
$heading.Range.InsertParagraphAfter()
$pA = $d.Paragraphs.Last
$pA.Range.InsertParagraphAfter()
$pB = $d.Paragraphs.Last
$pB.Range.InsertParagraphAfter()
$pC = $d.Paragraphs.Last
$pC.Range.InsertParagraphAfter()
$pD = $d.Paragraphs.Last
$pD.Range.InsertParagraphAfter()
$pE = $d.Paragraphs.Last

# --- Format & fill the heading paragraph ---
$hr = $heading.Range
$heading.SpaceAfter = 12
$heading.LineSpacingRule = 5
$hr.Font.Size = 12
$hr.Font.SizeBi = 12
$hr.InsertAfter("Variables de respuesta")

# --- Paragraph A: response-variable explanation ---
$pA.Range.InsertAfter("Las variables de respuesta que serán tomadas en cuenta en este experimento son el tiempo y la correctitud, ya que son la principal métrica para evaluar el desempeño de una implementación de Decision Tree. Las unidades utilizadas para la medición del tiempo serán los milisegundos, ya que proveen la suficiente precisión buscada para el desarrollo de las comparaciones.")

# --- Paragraph B stays blank ---

# --- Paragraph C: evaluation tool ---
$pC.Range.InsertAfter("Nuestra herramienta de evaluación para el tiempo será el método Stopwatch de C#.")

# --- Paragraph D stays blank ---

# --- Paragraph E: measurement methodology closing sentence ---
$pE.Range.InsertAfter("Debido a que la medición varía con esta herramienta, se probará cada caso 4 veces y se utilizará el tiempo promedio entre cada caso, logrando así una medida mas precisa del desempeño.")

"paragraphs=$($d.Paragraphs.Count)"
